$wb = $excel.ActiveWorkbook

# --- Sprint 2 sheet ---
$ws2 = $wb.Worksheets.Item("Sprint 2")

# Corrected daily "hours left" burndown numbers for sprint 2 (columns
# D:I = Day1..Day6, rows 2:5 = the 4 stories in this sprint).
$ws2.Cells.Item(2, 4).Value = 1
$ws2.Cells.Item(2, 5).Value = 1
$ws2.Cells.Item(2, 6).Value = 1
$ws2.Cells.Item(2, 7).Value = 1
$ws2.Cells.Item(2, 8).Value = 0
$ws2.Cells.Item(2, 9).Value = 0

$ws2.Cells.Item(3, 4).Value = 8
$ws2.Cells.Item(3, 5).Value = 8
$ws2.Cells.Item(3, 6).Value = 3
$ws2.Cells.Item(3, 7).Value = 3
$ws2.Cells.Item(3, 8).Value = 3
$ws2.Cells.Item(3, 9).Value = 0

$ws2.Cells.Item(4, 4).Value = 4
$ws2.Cells.Item(4, 5).Value = 4
$ws2.Cells.Item(4, 6).Value = 4
$ws2.Cells.Item(4, 7).Value = 2
$ws2.Cells.Item(4, 8).Value = 0
$ws2.Cells.Item(4, 9).Value = 0

$ws2.Cells.Item(5, 4).Value = 2
$ws2.Cells.Item(5, 5).Value = 2
$ws2.Cells.Item(5, 6).Value = 2
$ws2.Cells.Item(5, 7).Value = 2
$ws2.Cells.Item(5, 8).Value = 0
$ws2.Cells.Item(5, 9).Value = 0

# Drop the redundant "Remaining" helper column (J) - header text and the
# per-row formulas that duplicated the burndown calc.
$ws2.Range("J1:J9").ClearContents()

# "Actual Remaining Hours" row: sum each day's column directly instead of
# the previous day-over-day subtraction so a bad day1 entry can't throw
# off every later day's total.
$ws2.Range("D11").Formula = "=SUM(D2:D9)"
$ws2.Range("E11").Formula = "=SUM(E2:E9)"
$ws2.Range("F11").Formula = "=SUM(F2:F9)"
$ws2.Range("G11").Formula = "=SUM(G2:G9)"
$ws2.Range("H11").Formula = "=SUM(H2:H9)"
$ws2.Range("I11").Formula = "=SUM(I2:I9)"

# Show the "Sprint 2" tab (so time-left-for-sprint-2 is what's visible)
# with B3 selected.
$ws2.Activate()
$ws2.Range("B3").Select()
